$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): shift Date/Bays Opened/Total Opened/CFS one column right,
# and put the new "Day" header in column A. The old spacer cell (old E1) is dropped.
$ws.Range("E1").Value = $ws.Range("D1").Value2
$ws.Range("D1").Value = $ws.Range("C1").Value2
$ws.Range("C1").Value = $ws.Range("B1").Value2
$ws.Range("B1").Value = $ws.Range("A1").Value2
$ws.Range("A1").Value = "Day"
$ws.Range("A1").Style = "Normal"

# --- Data rows (2-24): column A becomes a plain sequential day number instead of
# the "Day N" text label; all other columns (B-G) stay as-is.
for ($i = 0; $i -lt 23; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
}

# --- Selection moves from I31 to A2
$ws.Range("A2").Select()
